# "Generate Report for Handoff"
#
# The localization-status report was regenerated: the zh-cn / de-de rows
# moved from "In Translation" to "Ready for handoff", and the two
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# were bumped forward to the new run's wall-clock time. Excel then
# auto-sized the (now wider) Status columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: File Name / Path / Extension / Publish URL / zh-cn / de-de / Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-30 16:47:19"

# --- zh-cn sheet: Status column (C) + Latest Handoff Datetime (H)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-30 16:47:09"

# --- de-de sheet: Status column (C) + Latest Handoff Datetime (H)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-30 16:47:19"

# The new "Ready for handoff" status string is wider than the old
# "In Translation" one, so Excel widens the Status column(s) on each
# sheet to fit it again.
$wsOverview.Columns("E:F").ColumnWidth = 16.3333333
$wsZhCn.Columns("C:C").ColumnWidth = 16.3333333
$wsDeDe.Columns("C:C").ColumnWidth = 16.3333333
